# "modified work study feature"
#
# The "Sheet" worksheet is a running log of work-study check-in/out
# entries (one per row, column A). This change appends a new batch of
# log entries for the 26th (14:51, 14:51 dup, 14:52, 14:53, 15:14,
# 15:28, 15:28 dup) right after the existing last entry (row 176).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$entries = @(
    "26 14:51>>> CFD893A460   Gary Tsai",
    "26 14:51>>> CFD893A460   Gary Tsai",
    "26 14:52>>> CFD893A460   Gary Tsai",
    "26 14:53>>> CFD893A460   Gary Tsai",
    "26 15:14>>> CFD893A460   Gary Tsai",
    "26 15:28>>> CFD893A460   Gary Tsai",
    "26 15:28>>> CFD893A460   Gary Tsai"
)

$startRow = 177
for ($i = 0; $i -lt $entries.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $entries[$i]
}
